# Auto-generated Excel COM-interop script to apply numeric corrections
# to the Sheets (Excalibur_Profits workbook) per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1112.0513  # H17
$ws.Cells.Item(17, 10).Value = 1115.8158  # J17
$ws.Cells.Item(17, 12).Value = 3347.4474  # L17
$ws.Cells.Item(17, 14).Value = -3683.4474  # N17

$ws.Cells.Item(70, 8).Value = 2524.2354  # H70
$ws.Cells.Item(70, 9).Value = 2450.6667  # I70
$ws.Cells.Item(70, 10).Value = 2540  # J70
$ws.Cells.Item(70, 11).Value = 7352.000100000001  # K70
$ws.Cells.Item(70, 12).Value = 7620  # L70
$ws.Cells.Item(70, 13).Value = -7082.000100000001  # M70
$ws.Cells.Item(70, 14).Value = -8160  # N70

$ws.Cells.Item(73, 8).Value = 2524.2354  # H73
$ws.Cells.Item(73, 9).Value = 2450.6667  # I73
$ws.Cells.Item(73, 10).Value = 2540  # J73
$ws.Cells.Item(73, 11).Value = 7352.000100000001  # K73
$ws.Cells.Item(73, 12).Value = 7620  # L73
$ws.Cells.Item(73, 13).Value = -6416.000100000001  # M73
$ws.Cells.Item(73, 14).Value = -9492  # N73

$ws.Cells.Item(76, 8).Value = 6363.1333  # H76
$ws.Cells.Item(76, 9).Value = 4490  # I76
$ws.Cells.Item(76, 10).Value = 7299.7  # J76
$ws.Cells.Item(76, 11).Value = 4490  # K76
$ws.Cells.Item(76, 12).Value = 7299.7  # L76
$ws.Cells.Item(76, 13).Value = -4175  # M76
$ws.Cells.Item(76, 14).Value = -7929.7  # N76

$ws.Cells.Item(79, 8).Value = 6363.1333  # H79
$ws.Cells.Item(79, 9).Value = 4490  # I79
$ws.Cells.Item(79, 10).Value = 7299.7  # J79
$ws.Cells.Item(79, 11).Value = 4490  # K79
$ws.Cells.Item(79, 12).Value = 7299.7  # L79
$ws.Cells.Item(79, 13).Value = -3398  # M79
$ws.Cells.Item(79, 14).Value = -9483.700000000001  # N79

$ws.Cells.Item(80, 8).Value = 4546518.5  # H80
$ws.Cells.Item(80, 9).Value = 11364408  # I80
$ws.Cells.Item(80, 10).Value = 1258.9166  # J80
$ws.Cells.Item(80, 11).Value = 34093224  # K80
$ws.Cells.Item(80, 12).Value = 3776.7498  # L80
$ws.Cells.Item(80, 13).Value = -34092226  # M80
$ws.Cells.Item(80, 14).Value = -5772.7498  # N80

$ws.Cells.Item(83, 8).Value = 4546518.5  # H83
$ws.Cells.Item(83, 9).Value = 11364408  # I83
$ws.Cells.Item(83, 10).Value = 1258.9166  # J83
$ws.Cells.Item(83, 11).Value = 102279672  # K83
$ws.Cells.Item(83, 12).Value = 11330.2494  # L83
$ws.Cells.Item(83, 13).Value = -102274680  # M83
$ws.Cells.Item(83, 14).Value = -21314.2494  # N83

$ws.Cells.Item(86, 8).Value = 4095.5356  # H86
$ws.Cells.Item(86, 9).Value = 2675.3845  # I86
$ws.Cells.Item(86, 11).Value = 2675.3845  # K86
$ws.Cells.Item(86, 13).Value = -1552.3845  # M86

$ws.Cells.Item(89, 8).Value = 4095.5356  # H89
$ws.Cells.Item(89, 9).Value = 2675.3845  # I89
$ws.Cells.Item(89, 11).Value = 13376.9225  # K89
$ws.Cells.Item(89, 13).Value = -7760.922500000001  # M89

$ws.Cells.Item(106, 8).Value = 3999.75  # H106
$ws.Cells.Item(106, 9).Value = 3999.75  # I106
$ws.Cells.Item(106, 11).Value = 3999.75  # K106
$ws.Cells.Item(106, 13).Value = -3368.75  # M106

$ws.Cells.Item(125, 8).Value = 1427.75  # H125
$ws.Cells.Item(125, 9).Value = 1144.5  # I125
$ws.Cells.Item(125, 11).Value = 10300.5  # K125
$ws.Cells.Item(125, 13).Value = -7840.5  # M125

$ws.Cells.Item(134, 8).Value = 70000  # H134
$ws.Cells.Item(134, 10).Value = 70000  # J134
$ws.Cells.Item(134, 12).Value = 70000  # L134
$ws.Cells.Item(134, 14).Value = -80140  # N134

$ws.Cells.Item(137, 8).Value = 7301.476  # H137
$ws.Cells.Item(137, 10).Value = 21198  # J137
$ws.Cells.Item(137, 12).Value = 63594  # L137
$ws.Cells.Item(137, 14).Value = -68694  # N137

$ws.Cells.Item(139, 8).Value = 85000  # H139
$ws.Cells.Item(139, 10).Value = 0  # J139
$ws.Cells.Item(139, 12).Value = 0  # L139
$ws.Cells.Item(139, 14).ClearContents()  # N139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(95, 8).Value = 62966.4  # H95
$ws.Cells.Item(95, 10).Value = 62966.4  # J95
$ws.Cells.Item(95, 12).Value = 62966.4  # L95
$ws.Cells.Item(95, 14).Value = -68458.39999999999  # N95

$ws.Cells.Item(102, 8).Value = 40662.332  # H102
$ws.Cells.Item(102, 9).Value = 59946  # I102
$ws.Cells.Item(102, 10).Value = 2095  # J102
$ws.Cells.Item(102, 11).Value = 59946  # K102
$ws.Cells.Item(102, 12).Value = 2095  # L102
$ws.Cells.Item(102, 13).Value = -58324  # M102
$ws.Cells.Item(102, 14).Value = -5339  # N102

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(51, 8).Value = 48500  # H51
$ws.Cells.Item(51, 10).Value = 48500  # J51
$ws.Cells.Item(51, 12).Value = 48500  # L51
$ws.Cells.Item(51, 14).Value = -49482  # N51

$ws.Cells.Item(92, 8).Value = 78517.336  # H92
$ws.Cells.Item(92, 10).Value = 78517.336  # J92
$ws.Cells.Item(92, 12).Value = 78517.336  # L92
$ws.Cells.Item(92, 14).Value = -83509.336  # N92

$ws.Cells.Item(95, 8).Value = 48762  # H95
$ws.Cells.Item(95, 10).Value = 48762  # J95
$ws.Cells.Item(95, 12).Value = 48762  # L95
$ws.Cells.Item(95, 14).Value = -54254  # N95

$ws.Cells.Item(119, 8).Value = 100761  # H119
$ws.Cells.Item(119, 10).Value = 100761  # J119
$ws.Cells.Item(119, 12).Value = 100761  # L119
$ws.Cells.Item(119, 14).Value = -110437  # N119

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 20680.74  # H31
$ws.Cells.Item(31, 9).Value = 8638.333000000001  # I31
$ws.Cells.Item(31, 11).Value = 8638.333000000001  # K31
$ws.Cells.Item(31, 13).Value = -8343.333000000001  # M31

$ws.Cells.Item(34, 8).Value = 20680.74  # H34
$ws.Cells.Item(34, 9).Value = 8638.333000000001  # I34
$ws.Cells.Item(34, 11).Value = 8638.333000000001  # K34
$ws.Cells.Item(34, 13).Value = -8436.333000000001  # M34

$ws.Cells.Item(52, 8).Value = 99989.5  # H52
$ws.Cells.Item(52, 9).Value = 0  # I52
$ws.Cells.Item(52, 10).Value = 99989.5  # J52
$ws.Cells.Item(52, 11).Value = 0  # K52
$ws.Cells.Item(52, 12).Value = 99989.5  # L52
$ws.Cells.Item(52, 13).ClearContents()  # M52
$ws.Cells.Item(52, 14).Value = -100577.5  # N52

$ws.Cells.Item(107, 8).Value = 614.6842  # H107
$ws.Cells.Item(107, 10).Value = 688.94116  # J107
$ws.Cells.Item(107, 12).Value = 688.94116  # L107
$ws.Cells.Item(107, 14).Value = -4528.94116  # N107

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 259.1579  # H33
$ws.Cells.Item(33, 10).Value = 898.5  # J33
$ws.Cells.Item(33, 12).Value = 5391  # L33
$ws.Cells.Item(33, 14).Value = -5957  # N33

$ws.Cells.Item(44, 8).Value = 2206.8333  # H44
$ws.Cells.Item(44, 9).Value = 995.8  # I44
$ws.Cells.Item(44, 10).Value = 3071.8572  # J44
$ws.Cells.Item(44, 11).Value = 2987.4  # K44
$ws.Cells.Item(44, 12).Value = 9215.571599999999  # L44
$ws.Cells.Item(44, 13).Value = -2589.4  # M44
$ws.Cells.Item(44, 14).Value = -10011.5716  # N44

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 880648.1  # H126
$ws.Cells.Item(126, 9).Value = 1284793.5  # I126
$ws.Cells.Item(126, 10).Value = 4999.8335  # J126
$ws.Cells.Item(126, 11).Value = 3854380.5  # K126
$ws.Cells.Item(126, 12).Value = 14999.5005  # L126
$ws.Cells.Item(126, 13).Value = -3851910.5  # M126
$ws.Cells.Item(126, 14).Value = -19939.5005  # N126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1500  # H16
$ws.Cells.Item(16, 9).Value = 0  # I16
$ws.Cells.Item(16, 10).Value = 1500  # J16
$ws.Cells.Item(16, 11).Value = 0  # K16
$ws.Cells.Item(16, 12).Value = 1500  # L16
$ws.Cells.Item(16, 13).ClearContents()  # M16
$ws.Cells.Item(16, 14).Value = -1840  # N16

$ws.Cells.Item(61, 8).Value = 2718.8  # H61
$ws.Cells.Item(61, 9).Value = 2348.3  # I61
$ws.Cells.Item(61, 10).Value = 3459.8  # J61
$ws.Cells.Item(61, 11).Value = 2348.3  # K61
$ws.Cells.Item(61, 12).Value = 3459.8  # L61
$ws.Cells.Item(61, 13).Value = -2146.3  # M61
$ws.Cells.Item(61, 14).Value = -3863.8  # N61

$ws.Cells.Item(113, 8).Value = 2718.8  # H113
$ws.Cells.Item(113, 9).Value = 2348.3  # I113
$ws.Cells.Item(113, 10).Value = 3459.8  # J113
$ws.Cells.Item(113, 11).Value = 2348.3  # K113
$ws.Cells.Item(113, 12).Value = 3459.8  # L113
$ws.Cells.Item(113, 13).Value = -178.3000000000002  # M113
$ws.Cells.Item(113, 14).Value = -7799.8  # N113

$ws.Cells.Item(132, 8).Value = 38518.145  # H132
$ws.Cells.Item(132, 9).Value = 50925.4  # I132
$ws.Cells.Item(132, 10).Value = 7500  # J132
$ws.Cells.Item(132, 11).Value = 152776.2  # K132
$ws.Cells.Item(132, 12).Value = 22500  # L132
$ws.Cells.Item(132, 13).Value = -150246.2  # M132
$ws.Cells.Item(132, 14).Value = -27560  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 80333.836  # H62
$ws.Cells.Item(62, 9).Value = 302333.34  # I62
$ws.Cells.Item(62, 10).Value = 6334  # J62
$ws.Cells.Item(62, 11).Value = 302333.34  # K62
$ws.Cells.Item(62, 12).Value = 6334  # L62
$ws.Cells.Item(62, 13).Value = -301709.34  # M62
$ws.Cells.Item(62, 14).Value = -7582  # N62

$ws.Cells.Item(65, 8).Value = 80333.836  # H65
$ws.Cells.Item(65, 9).Value = 302333.34  # I65
$ws.Cells.Item(65, 10).Value = 6334  # J65
$ws.Cells.Item(65, 11).Value = 1511666.7  # K65
$ws.Cells.Item(65, 12).Value = 31670  # L65
$ws.Cells.Item(65, 13).Value = -1508546.7  # M65
$ws.Cells.Item(65, 14).Value = -37910  # N65

$ws.Cells.Item(107, 8).Value = 1254.9584  # H107
$ws.Cells.Item(107, 9).Value = 642.2273  # I107
$ws.Cells.Item(107, 11).Value = 1926.6819  # K107
$ws.Cells.Item(107, 13).Value = -6.681900000000041  # M107

$ws.Cells.Item(116, 8).Value = 200666.5  # H116
$ws.Cells.Item(116, 10).Value = 200666.5  # J116
$ws.Cells.Item(116, 12).Value = 200666.5  # L116
$ws.Cells.Item(116, 14).Value = -209844.5  # N116

$ws.Cells.Item(126, 8).Value = 1772.421  # H126
$ws.Cells.Item(126, 9).Value = 1648.9445  # I126
$ws.Cells.Item(126, 11).Value = 4946.833500000001  # K126
$ws.Cells.Item(126, 13).Value = -2476.833500000001  # M126

$ws.Cells.Item(132, 8).Value = 3270532.2  # H132
$ws.Cells.Item(132, 9).Value = 3970631.8  # I132
$ws.Cells.Item(132, 10).Value = 3401  # J132
$ws.Cells.Item(132, 11).Value = 11911895.4  # K132
$ws.Cells.Item(132, 12).Value = 10203  # L132
$ws.Cells.Item(132, 13).Value = -11909365.4  # M132
$ws.Cells.Item(132, 14).Value = -15263  # N132
